# Crypto_Fall_14_Group.docx edit script
# Implements: added Crypto++ include list, PIN values, "$"1,000,000 and packet
# padding figures, "50" connection limit, and assorted typo fixes.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert the block of #include lines (plus trailing blank line) right
#    before the "Running the Program" heading.  The insertion point is
#    chosen at the end of the *previous* plain paragraph so the new runs
#    pick up normal (non-heading) character formatting.
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "specific libraries that we use are:"
$found = $find.Execute()
if ($found) {
    $r = $find.Parent
    $r.Collapse(0)
    $includeText = "`r" + '#include "includes/cryptopp/sha.h"' + "`r" + `
                   '#include "includes/cryptopp/hex.h"' + "`r" + `
                   '#include "includes/cryptopp/aes.h"' + "`r" + `
                   '#include "includes/cryptopp/ccm.h"' + "`r" + `
                   '#include "includes/cryptopp/gcm.h"' + "`r" + `
                   '#include "includes/cryptopp/osrng.h"'
    $r.InsertBefore($includeText)
}

# Append the trailing blank paragraph after the last #include line.
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Text = 'osrng.h"'
$found2 = $find2.Execute()
if ($found2) {
    $r2 = $find2.Parent
    $r2.Collapse(0)
    $r2.InsertParagraphAfter()
}

# Apply "space after = 0" to the 6 include paragraphs plus the blank line
# that follows them (7 paragraphs total).
$idx = 0
$startIdx = -1
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -like "*#include*sha.h*") {
        $startIdx = $idx
    }
}
if ($startIdx -gt 0) {
    for ($i = $startIdx; $i -le ($startIdx + 6); $i++) {
        $d.Paragraphs($i).Range.ParagraphFormat.SpaceAfter = 0
    }
}

# ---------------------------------------------------------------------------
# 2. Fix "program take one" -> "program takes one" typo.
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("program take one command", $true, $false, $false, $false, $false, `
               $true, 1, $false, "program takes one command", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Fill in the PIN column of the accounts table.
# ---------------------------------------------------------------------------
$t = $d.Tables(1)
$t.Cell(2, 2).Range.Text = "123456"
$t.Cell(3, 2).Range.Text = "456789"
$t.Cell(4, 2).Range.Text = "654321"

# ---------------------------------------------------------------------------
# 4. Remove the duplicated space "Withdraw,  Transfer" -> "Withdraw, Transfer"
#    and delete the stray blank paragraph that followed the transactions
#    paragraph.
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("Withdraw,  Transfer", $true, $false, $false, $false, $false, `
               $true, 1, $false, "Withdraw, Transfer", 2) | Out-Null

$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -like "*Deposit and Balance*") {
        $d.Paragraphs($idx + 1).Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 5. Balance the parenthesis typo: "Nbank2)" -> "Nbank2))"
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("bank2)", $true, $false, $false, $false, $false, `
               $true, 1, $false, "bank2))", 2) | Out-Null

# ---------------------------------------------------------------------------
# 6. "accept XX connections" -> "accept 50 connections" (remove highlight).
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "accept XX connections"
$found = $find.Execute()
if ($found) {
    $r = $find.Parent
    $r.Text = "accept 50 connections"
    $r.HighlightColorIndex = 0
}

# ---------------------------------------------------------------------------
# 7. Maximum account amount + packet padding paragraph rewrite.
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "XXXXXXXXXXXXX"
$found = $find.Execute()
if ($found) {
    $r = $find.Parent
    $r.Text = "`$1,000,000"
    $r.HighlightColorIndex = 0
}

$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "."
$find.Text = "is `$1,000,000."
$found = $find.Execute()
if ($found) {
    $r = $find.Parent
    $r.Collapse(0)
    $r.InsertAfter("  ")
}

$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "All Packets are padded to XXXX bytes to ensure no information leakage."
$found = $find.Execute()
if ($found) {
    $r = $find.Parent
    $r.Text = "All Packets are padded to 984 bytes to ensure no information leakage.  984 was selected as it provides the largest encrypted block smaller than 1024 in AES.  The packets will be sent at random intervals not to exceed 1 second in order to help thwart timing attacks."
    $r.HighlightColorIndex = 0
}
